$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new date columns before column C (which holds the "Jun_10"
# history). This shifts the existing column C (data + formatting) two
# places to the right, landing it in column E.
$ws.Columns("C:D").Insert()

# Re-populate the header row (row 1) with the new date labels.
# Column C is set before column B so the shared-string table records
# "Jun_15" ahead of "Jun_17" (matching first-use append order).
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"
$ws.Range("D1").Value = "Jun_13"

# The two freshly inserted columns start out empty; every other column in
# this table uses the "UN" placeholder for rows with no specific rating
# change, so mirror that for the new C and D columns across all data rows.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Cosmetic: keep the 8-character custom width that column C originally
# had, now applied across the three date columns C, D and E.
$ws.Columns("C").ColumnWidth = 7.166666666666667
$ws.Columns("D").ColumnWidth = 7.166666666666667
$ws.Columns("E").ColumnWidth = 7.166666666666667
